# Updates the cryptos price/volume table to the latest scraped snapshot
# (GitHub Actions refresh run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin / Link column swaps (three coins re-ranked this run) ---
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"

# --- Price column (D). Force text so values like "29.085.66" or
# "0.00000000126" are stored verbatim instead of being parsed as
# numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "29.085.66"
$ws.Cells.Item(3, 4).Value = "1.834.26"
$ws.Cells.Item(4, 4).Value = "1.001"
$ws.Cells.Item(5, 4).Value = "243.92"
$ws.Cells.Item(6, 4).Value = "0.6292"
$ws.Cells.Item(7, 4).Value = "1.002"
$ws.Cells.Item(8, 4).Value = "0.07503"
$ws.Cells.Item(9, 4).Value = "0.2927"
$ws.Cells.Item(10, 4).Value = "23.20"
$ws.Cells.Item(11, 4).Value = "0.07722"
$ws.Cells.Item(12, 4).Value = "1.838.68"
$ws.Cells.Item(13, 4).Value = "4.997"
$ws.Cells.Item(14, 4).Value = "0.6691"
$ws.Cells.Item(15, 4).Value = "82.66"
$ws.Cells.Item(16, 4).Value = "0.000009342"
$ws.Cells.Item(17, 4).Value = "6.003"
$ws.Cells.Item(18, 4).Value = "29.121.10"
$ws.Cells.Item(19, 4).Value = "2.081.01"
$ws.Cells.Item(20, 4).Value = "12.59"
$ws.Cells.Item(21, 4).Value = "223.70"
$ws.Cells.Item(23, 4).Value = "7.130"
$ws.Cells.Item(24, 4).Value = "1.002"
$ws.Cells.Item(25, 4).Value = "159.58"
$ws.Cells.Item(26, 4).Value = "0.1403"
$ws.Cells.Item(27, 4).Value = "8.512"
$ws.Cells.Item(28, 4).Value = "17.96"
$ws.Cells.Item(29, 4).Value = "1.500"
$ws.Cells.Item(30, 4).Value = "0.05768"
$ws.Cells.Item(31, 4).Value = "4.158"
$ws.Cells.Item(32, 4).Value = "4.065"
$ws.Cells.Item(33, 4).Value = "1.205"
$ws.Cells.Item(34, 4).Value = "0.7511"
$ws.Cells.Item(35, 4).Value = "1.848"
$ws.Cells.Item(36, 4).Value = "1.139"
$ws.Cells.Item(37, 4).Value = "2.673"
$ws.Cells.Item(38, 4).Value = "2.766"
$ws.Cells.Item(39, 4).Value = "1.222.54"
$ws.Cells.Item(40, 4).Value = "0.01786"
$ws.Cells.Item(41, 4).Value = "6.565"
$ws.Cells.Item(42, 4).Value = "0.8935"
$ws.Cells.Item(44, 4).Value = "102.18"
$ws.Cells.Item(45, 4).Value = "0.00000000126"
$ws.Cells.Item(46, 4).Value = "1.981.55"
$ws.Cells.Item(47, 4).Value = "65.74"
$ws.Cells.Item(48, 4).Value = "0.07835"
$ws.Cells.Item(49, 4).Value = "0.5092"
$ws.Cells.Item(50, 4).Value = "0.4072"
$ws.Cells.Item(51, 4).Value = "9.023"

# --- Volume(1h) column (E) ---
$ws.Cells.Item(2, 5).Value = "  +0.21%  "
$ws.Cells.Item(3, 5).Value = "  +0.15%  "
$ws.Cells.Item(4, 5).Value = "  +0.15%  "
$ws.Cells.Item(5, 5).Value = "  +0.65%  "
$ws.Cells.Item(6, 5).Value = "  +0.65%  "
$ws.Cells.Item(7, 5).Value = "  +0.22%  "
$ws.Cells.Item(8, 5).Value = "  -0.70%  "
$ws.Cells.Item(9, 5).Value = "  +0.26%  "
$ws.Cells.Item(10, 5).Value = "  +2.92%  "
$ws.Cells.Item(11, 5).Value = "  -0.07%  "
$ws.Cells.Item(12, 5).Value = "  +0.35%  "
$ws.Cells.Item(13, 5).Value = "  +1.03%  "
$ws.Cells.Item(14, 5).Value = "  +0.67%  "
$ws.Cells.Item(15, 5).Value = "  -0.06%  "
$ws.Cells.Item(16, 5).Value = "  -6.99%  "
$ws.Cells.Item(17, 5).Value = "  -0.46%  "
$ws.Cells.Item(18, 5).Value = "  +0.46%  "
$ws.Cells.Item(19, 5).Value = "  +0.00%  "
$ws.Cells.Item(20, 5).Value = "  +2.17%  "
$ws.Cells.Item(21, 5).Value = "  -1.39%  "
$ws.Cells.Item(22, 5).Value = "  +0.37%  "
$ws.Cells.Item(23, 5).Value = "  -0.42%  "
$ws.Cells.Item(24, 5).Value = "  +0.19%  "
$ws.Cells.Item(25, 5).Value = "  +0.96%  "
$ws.Cells.Item(26, 5).Value = "  +2.19%  "
$ws.Cells.Item(27, 5).Value = "  +0.54%  "
$ws.Cells.Item(28, 5).Value = "  +0.14%  "
$ws.Cells.Item(29, 5).Value = "  +0.77%  "
$ws.Cells.Item(30, 5).Value = "  +10.97%  "
$ws.Cells.Item(31, 5).Value = "  +1.69%  "
$ws.Cells.Item(32, 5).Value = "  +1.17%  "
$ws.Cells.Item(33, 5).Value = "  +0.79%  "
$ws.Cells.Item(34, 5).Value = "  +1.97%  "
$ws.Cells.Item(35, 5).Value = "  +0.08%  "
$ws.Cells.Item(36, 5).Value = "  -0.04%  "
$ws.Cells.Item(37, 5).Value = "  -0.87%  "
$ws.Cells.Item(38, 5).Value = "  +0.19%  "
$ws.Cells.Item(39, 5).Value = "  -1.88%  "
$ws.Cells.Item(40, 5).Value = "  +0.10%  "
$ws.Cells.Item(41, 5).Value = "  +3.63%  "
$ws.Cells.Item(42, 5).Value = "  -0.24%  "
$ws.Cells.Item(43, 5).Value = "  +0.24%  "
$ws.Cells.Item(44, 5).Value = "  +0.82%  "
$ws.Cells.Item(45, 5).Value = "  +2.52%  "
$ws.Cells.Item(46, 5).Value = "  +0.12%  "
$ws.Cells.Item(47, 5).Value = "  +2.67%  "
$ws.Cells.Item(48, 5).Value = "  +13.35%  "
$ws.Cells.Item(49, 5).Value = "  -0.36%  "
$ws.Cells.Item(50, 5).Value = "  +0.81%  "
$ws.Cells.Item(51, 5).Value = "  +1.90%  "

